$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.791111874039018
$ws.Range("C2").Value = 0.03324406318822071
$ws.Range("E2").Value = 0.0642612727440639
$ws.Range("F2").Value = 3.87690584660406
$ws.Range("G2").Value = 0.002637231989988654
$ws.Range("J2").Value = 0.2446161763393331
$ws.Range("K2").Value = 1.28578184861567
$ws.Range("L2").Value = 0.2276476898771023
$ws.Range("M2").Value = 0.373626407520959
$ws.Range("N2").Value = 4.426523778078547
$ws.Range("B3").Value = 1.757513110745663
$ws.Range("C3").Value = 0.02894370036712246
$ws.Range("E3").Value = 0.06444610348159152
$ws.Range("F3").Value = 3.862069570967918
$ws.Range("G3").Value = 0.002641523552572135
$ws.Range("J3").Value = 0.2452501559365938
$ws.Range("K3").Value = 1.249973112376694
$ws.Range("L3").Value = 0.2264756928067584
$ws.Range("M3").Value = 0.3683393834281361
$ws.Range("N3").Value = 4.433308292482565
$ws.Range("B4").Value = 1.737860886305128
$ws.Range("C4").Value = 0.02630019188752897
$ws.Range("E4").Value = 0.0645762509700889
$ws.Range("F4").Value = 3.854390642668804
$ws.Range("G4").Value = 0.002644298649288543
$ws.Range("J4").Value = 0.2456539588095357
$ws.Range("K4").Value = 1.228737739201136
$ws.Range("L4").Value = 0.2258642718633226
$ws.Range("M4").Value = 0.3652911447098575
$ws.Range("N4").Value = 4.438317727532805
$ws.Range("B5").Value = 1.730098405744656
$ws.Range("C5").Value = 0.02522208149501637
$ws.Range("E5").Value = 0.0646334926623009
$ws.Range("F5").Value = 3.851621323845677
$ws.Range("G5").Value = 0.002645464855429003
$ws.Range("J5").Value = 0.2458221697563916
$ws.Range("K5").Value = 1.220273009258733
$ws.Range("L5").Value = 0.2256423957680127
$ws.Range("M5").Value = 0.3640988227680637
$ws.Range("N5").Value = 4.440571198703964
$ws.Range("B6").Value = 1.728824311699526
$ws.Range("C6").Value = 0.02504300763838785
$ws.Range("E6").Value = 0.06464325208021204
$ws.Range("F6").Value = 3.851183223205481
$ws.Range("G6").Value = 0.002645660640527362
$ws.Range("J6").Value = 0.245850322067831
$ws.Range("K6").Value = 1.218878851881442
$ws.Range("L6").Value = 0.2256072035261738
$ws.Range("M6").Value = 0.3639038520627906
$ws.Range("N6").Value = 4.440958193250836
$ws.Range("B7").Value = 1.73775520269686
$ws.Range("C7").Value = 0.02628565570772423
$ws.Range("E7").Value = 0.06457700590181403
$ws.Range("F7").Value = 3.854351837203637
$ws.Range("G7").Value = 0.00264431423384233
$ws.Range("J7").Value = 0.2456562125476864
$ws.Range("K7").Value = 1.228622816311571
$ws.Range("L7").Value = 0.22586116899312
$ws.Range("M7").Value = 0.365274862667448
$ws.Range("N7").Value = 4.438347260025594
$ws.Range("B8").Value = 1.779324200381268
$ws.Range("C8").Value = 0.03176188088598053
$ws.Range("E8").Value = 0.06432155402613216
$ws.Range("F8").Value = 3.871493464625843
$ws.Range("G8").Value = 0.002638682719411121
$ws.Range("J8").Value = 0.2448317569989413
$ws.Range("K8").Value = 1.273278996583144
$ws.Range("L8").Value = 0.2272211759337708
$ws.Range("M8").Value = 0.3717623809831672
$ws.Range("N8").Value = 4.428687917841529
$ws.Range("B9").Value = 1.868597799042334
$ws.Range("C9").Value = 0.04248083778131218
$ws.Range("E9").Value = 0.06395212277945728
$ws.Range("F9").Value = 3.916457569267337
$ws.Range("G9").Value = 0.002628745497978854
$ws.Range("J9").Value = 0.2433301969132131
$ws.Range("K9").Value = 1.366821023761275
$ws.Range("L9").Value = 0.2307435917497074
$ws.Range("M9").Value = 0.3860537780634488
$ws.Range("N9").Value = 4.416445182425377
$ws.Range("B10").Value = 1.938926891603842
$ws.Range("C10").Value = 0.05035077170568059
$ws.Range("E10").Value = 0.06375998093214186
$ws.Range("F10").Value = 3.956417631701242
$ws.Range("G10").Value = 0.00262211174993811
$ws.Range("J10").Value = 0.242296975610385
$ws.Range("K10").Value = 1.439209906518983
$ws.Range("L10").Value = 0.233849609326839
$ws.Range("M10").Value = 0.3975098029892905
$ws.Range("N10").Value = 4.411542976555268
$ws.Range("B11").Value = 1.971953734095564
$ws.Range("C11").Value = 0.05393138269511155
$ws.Range("E11").Value = 0.06368960051483441
$ws.Range("F11").Value = 3.976102141404112
$ws.Range("G11").Value = 0.00261923720935231
$ws.Range("J11").Value = 0.241842078118454
$ws.Range("K11").Value = 1.472942681550307
$ws.Range("L11").Value = 0.2353744447371682
$ws.Range("M11").Value = 0.4029290688621856
$ws.Range("N11").Value = 4.410203539660969
$ws.Range("B12").Value = 1.984608854769363
$ws.Range("C12").Value = 0.05528745238316901
$ws.Range("E12").Value = 0.06366538223589124
$ws.Range("F12").Value = 3.983772793668692
$ws.Range("G12").Value = 0.002618169168437735
$ws.Range("J12").Value = 0.2416719921881061
$ws.Range("K12").Value = 1.485832108819977
$ws.Range("L12").Value = 0.2359678833457721
$ws.Range("M12").Value = 0.4050110607444353
$ws.Range("N12").Value = 4.40982454443909
$ws.Range("B13").Value = 1.981876741586632
$ws.Range("C13").Value = 0.05499538947078975
$ws.Range("E13").Value = 0.06367049009381009
$ws.Range("F13").Value = 3.982111151188747
$ws.Range("G13").Value = 0.002618398280665207
$ws.Range("J13").Value = 0.2417085265561809
$ws.Range("K13").Value = 1.483050998461977
$ws.Range("L13").Value = 0.2358393645474237
$ws.Range("M13").Value = 0.4045613405793347
$ws.Range("N13").Value = 4.409900462711931
$ws.Range("B14").Value = 1.972991901341459
$ws.Range("C14").Value = 0.05404294342524452
$ws.Range("E14").Value = 0.06368755938187043
$ws.Range("F14").Value = 3.97672887088504
$ws.Range("G14").Value = 0.002619148930909389
$ws.Range("J14").Value = 0.2418280414803462
$ws.Range("K14").Value = 1.474000785782806
$ws.Range("L14").Value = 0.2354229467117932
$ws.Range("M14").Value = 0.4030997580005646
$ws.Range("N14").Value = 4.410169788876317
$ws.Range("B15").Value = 1.967569026239516
$ws.Range("C15").Value = 0.05345956761935611
$ws.Range("E15").Value = 0.0636983312521231
$ws.Range("F15").Value = 3.973460268645738
$ws.Range("G15").Value = 0.002619611390367774
$ws.Range("J15").Value = 0.2419015309177066
$ws.Range("K15").Value = 1.468472322130594
$ws.Range("L15").Value = 0.2351699624193486
$ws.Range("M15").Value = 0.4022083800064706
$ws.Range("N15").Value = 4.41035146163064
$ws.Range("B16").Value = 1.936789304756473
$ws.Range("C16").Value = 0.0501167890926979
$ws.Range("E16").Value = 0.06376492177248672
$ws.Range("F16").Value = 3.955161510395641
$ws.Range("G16").Value = 0.002622302480197557
$ws.Range("J16").Value = 0.2423270084496068
$ws.Range("K16").Value = 1.437021553761184
$ws.Range("L16").Value = 0.2337522030145749
$ws.Range("M16").Value = 0.397159819054778
$ws.Range("N16").Value = 4.411648446066735
$ws.Range("B17").Value = 1.918171657092614
$ws.Range("C17").Value = 0.04806630635229681
$ws.Range("E17").Value = 0.06381012288097843
$ws.Range("F17").Value = 3.944321605067003
$ws.Range("G17").Value = 0.002623989975126964
$ws.Range("J17").Value = 0.2425918974173467
$ws.Range("K17").Value = 1.417933208684957
$ws.Range("L17").Value = 0.2329110594118617
$ws.Range("M17").Value = 0.3941158881535287
$ws.Range("N17").Value = 4.412672315878837
$ws.Range("B18").Value = 1.907560587998717
$ws.Range("C18").Value = 0.04688696803491155
$ws.Range("E18").Value = 0.06383772480824934
$ws.Range("F18").Value = 3.938228579002782
$ws.Range("G18").Value = 0.002624974060808464
$ws.Range("J18").Value = 0.2427456777274948
$ws.Range("K18").Value = 1.407029644147144
$ws.Range("L18").Value = 0.2324377920599474
$ws.Range("M18").Value = 0.392384667949834
$ws.Range("N18").Value = 4.413345036918599
$ws.Range("B19").Value = 1.903984578014956
$ws.Range("C19").Value = 0.04648766928423242
$ws.Range("E19").Value = 0.06384734620077559
$ws.Range("F19").Value = 3.936189941833646
$ws.Range("G19").Value = 0.002625309574463201
$ws.Range("J19").Value = 0.2427979895453198
$ws.Range("K19").Value = 1.403350855937646
$ws.Range("L19").Value = 0.2322793637170051
$ws.Range("M19").Value = 0.3918018691704859
$ws.Range("N19").Value = 4.413587199949717
$ws.Range("B20").Value = 1.920143466967033
$ws.Range("C20").Value = 0.0482845784009811
$ws.Range("E20").Value = 0.06380514530446924
$ws.Range("F20").Value = 3.94546085624782
$ws.Range("G20").Value = 0.002623808943804834
$ws.Range("J20").Value = 0.2425635522453646
$ws.Range("K20").Value = 1.419957375907359
$ws.Range("L20").Value = 0.2329995106616565
$ws.Range("M20").Value = 0.3944378949638292
$ws.Range("N20").Value = 4.412554647275343
$ws.Range("B21").Value = 1.975597562725852
$ws.Range("C21").Value = 0.05432269456206029
$ws.Range("E21").Value = 0.06368247979943042
$ws.Range("F21").Value = 3.97830390037555
$ws.Range("G21").Value = 0.00261892789155843
$ws.Range("J21").Value = 0.2417928780531486
$ws.Range("K21").Value = 1.476655915370117
$ws.Range("L21").Value = 0.2355448246225791
$ws.Range("M21").Value = 0.4035282510072307
$ws.Range("N21").Value = 4.410087200143835
$ws.Range("B22").Value = 2.012705978393853
$ws.Range("C22").Value = 0.05826998691551921
$ws.Range("E22").Value = 0.06361648682126653
$ws.Range("F22").Value = 4.001030950855721
$ws.Range("G22").Value = 0.002615857198327332
$ws.Range("J22").Value = 0.2413018689057811
$ws.Range("K22").Value = 1.514385310518406
$ws.Range("L22").Value = 0.2373016573888407
$ws.Range("M22").Value = 0.4096431912468788
$ws.Range("N22").Value = 4.40922202188878
$ws.Range("B23").Value = 1.992821309437261
$ws.Range("C23").Value = 0.05616311729993129
$ws.Range("E23").Value = 0.06365041633690538
$ws.Range("F23").Value = 3.98878562553756
$ws.Range("G23").Value = 0.002617485198363172
$ws.Range("J23").Value = 0.2415627704048049
$ws.Range("K23").Value = 1.494186741461903
$ws.Range("L23").Value = 0.2363554877334479
$ws.Range("M23").Value = 0.4063636421615371
$ws.Range("N23").Value = 4.409615341233973
$ws.Range("B24").Value = 1.919251723702757
$ws.Range("C24").Value = 0.04818589913796245
$ws.Range("E24").Value = 0.06380739063482732
$ws.Range("F24").Value = 3.944945367872478
$ws.Range("G24").Value = 0.002623890744652448
$ws.Range("J24").Value = 0.2425763624492117
$ws.Range("K24").Value = 1.419042030072688
$ws.Range("L24").Value = 0.2329594897097707
$ws.Range("M24").Value = 0.394292257195751
$ws.Range("N24").Value = 4.412607583333511
$ws.Range("B25").Value = 1.843615388722611
$ws.Range("C25").Value = 0.03958250084130555
$ws.Range("E25").Value = 0.06403807758589153
$ws.Range("F25").Value = 3.90307795582865
$ws.Range("G25").Value = 0.002631316106071092
$ws.Range("J25").Value = 0.2437241029852562
$ws.Range("K25").Value = 1.340873771207811
$ws.Range("L25").Value = 0.2296994577567091
$ws.Range("M25").Value = 0.382019629976746
$ws.Range("N25").Value = 4.419039106200003
